$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 769 (shifts existing rows 769:810 down to 770:811)
$ws.Rows.Item(769).Insert()

# Populate the newly inserted row with the new data point.
# Force column A to stay plain text (matches the rest of the sheet) so
# Excel doesn't auto-convert the "yyyy/mm/dd" looking string into a date
# serial, then reset the style so no stray number-format sticks to the cell.
$ws.Range("A769").NumberFormat = "@"
$ws.Range("A769").Value = "2026/02/02"
$ws.Range("A769").Style = "Normal"
$ws.Range("B769").Value = "月"
$ws.Range("C769").Value = 16
$ws.Range("D769").Value = 185
